$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 199
$ws.Range("B4").Value = 144
$ws.Range("C4").Value = 44
$ws.Range("D4").Value = 11

$ws.Range("A10").Value = 462
$ws.Range("B10").Value = 341
$ws.Range("C10").Value = 49
$ws.Range("D10").Value = 59

$ws.Range("C67").Value = 0.0000035451000000000002
$ws.Range("D67").Value = 0.034616000000000001
$ws.Range("E67").Value = 1
$ws.Range("F67").Value = 0.19928000000000001
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 0.45846999999999999
$ws.Range("I67").Value = 0.53147999999999995
$ws.Range("J67").Value = 1
$ws.Range("K67").Value = 1
$ws.Range("L67").Value = 0.023508000000000001
$ws.Range("M67").Value = 0.000099125999999999999
$ws.Range("N67").Value = 1
$ws.Range("O67").Value = 0.019701
$ws.Range("P67").Value = 0.00036822999999999999
$ws.Range("Q67").Value = 1
$ws.Range("C68").Value = 0.39478999999999997
$ws.Range("D68").Value = 0.42201
$ws.Range("E68").Value = 0.49748999999999999
$ws.Range("F68").Value = 0.071894
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 0.0010757
$ws.Range("I68").Value = 0.088457999999999995
$ws.Range("J68").Value = 0.22681000000000001
$ws.Range("K68").Value = 0.63590000000000002
$ws.Range("L68").Value = 0.00232
$ws.Range("M68").Value = 0.52734999999999999
$ws.Range("N68").Value = 1
$ws.Range("O68").Value = 0.033514000000000002
$ws.Range("P68").Value = 0.17765
$ws.Range("Q68").Value = 1
$ws.Range("C69").Value = 0.21387
$ws.Range("D69").Value = 0.069061999999999998
$ws.Range("E69").Value = 1
$ws.Range("F69").Value = 0.19481999999999999
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 0.20172000000000001
$ws.Range("I69").Value = 0.087661000000000003
$ws.Range("J69").Value = 0.12146999999999999
$ws.Range("K69").Value = 0.19097
$ws.Range("L69").Value = 1
$ws.Range("C70").Value = 1
$ws.Range("D70").Value = 0.011858
$ws.Range("E70").Value = 1
$ws.Range("F70").Value = 1
$ws.Range("G70").Value = 0.025447000000000001
$ws.Range("H70").Value = 0.33783000000000002
$ws.Range("I70").Value = 0.52734999999999999
$ws.Range("J70").Value = 0.124
$ws.Range("K70").Value = 0.17765
$ws.Range("L70").Value = 1
$ws.Range("M70").Value = 0.00000029058
$ws.Range("N70").Value = 1
$ws.Range("O70").Value = 1
$ws.Range("P70").Value = 0.0073923000000000001
$ws.Range("Q70").Value = 1
$ws.Range("D71").Value = 1
$ws.Range("H71").Value = 0.52932000000000001
$ws.Range("C72").Value = 0.20069999999999999
$ws.Range("D72").Value = 0.12712000000000001
$ws.Range("E72").Value = 0.073043999999999998
$ws.Range("F72").Value = 0.097156000000000006
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 0.0071541
$ws.Range("I72").Value = 0.050233
$ws.Range("J72").Value = 0.011205
$ws.Range("K72").Value = 1
$ws.Range("L72").Value = 0.0052781
$ws.Range("M72").Value = 0.000054635999999999997
$ws.Range("N72").Value = 0.41596
$ws.Range("O72").Value = 0.042472999999999997
$ws.Range("P72").Value = 0.0074916999999999996
$ws.Range("Q72").Value = 1
